$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8309265375137329
$ws.Range("B1").Value = 2.858924150466919
$ws.Range("C1").Value = 4.694247722625732
$ws.Range("D1").Value = 2.810256242752075
$ws.Range("E1").Value = 1.447497725486755
